$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.840.05"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.48%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.644.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +0.42%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.98"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.26%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.05%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.33%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.30%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.875.32"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.59%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.633.15"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.10%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.29"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.876.56"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.35"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +3.35%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +7.68%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.44"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +6.02%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.00"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.38%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.67%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +0.65%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.87%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0511"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.58%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.18"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.00%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.38"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.85%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.44%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.67%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.80%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.247.14"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.27%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.533"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +1.16%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.70%  "
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.31%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.80%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.70%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.788.53"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -4.54%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "60.98"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +1.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.65"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.66%  "
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +6.43%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +1.19%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.52"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.40%  "
